$wb = $excel.ActiveWorkbook

$wsContainers = $wb.Worksheets.Item("Containers")
$wsDetails    = $wb.Worksheets.Item("Details")

# --- Details sheet: insert a new column before F for a "Comments" field ---
$wsDetails.Columns.Item(6).Insert()

# Header cell for the new column, formatted like its neighbouring headers
$wsDetails.Range("F4").Value = "Comments"

# Widen the new Comments column so the header/content is readable
$wsDetails.Columns.Item(6).ColumnWidth = 28.5

# --- Update saved selections / active sheet to match latest user interaction ---
# Containers sheet was last clicked on cell J13 ...
[void]$wsContainers.Activate()
[void]$wsContainers.Range("J13").Select()

# ... then the user finished on the Details sheet at cell F8, which is
# the sheet left active/selected when the workbook was saved.
[void]$wsDetails.Activate()
[void]$wsDetails.Range("F8").Select()
